$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("土地" / land) ---
# Remove the single data row, keep header row, and add two new header
# columns: "portion" and "total" (area * portion), per issue #5.
$s1 = $wb.Worksheets.Item(1)
$s1.Rows.Item(2).Delete()
$s1.Range("O1").Copy($s1.Range("P1"))
$s1.Range("P1").Value = "portion"
$s1.Range("O1").Copy($s1.Range("Q1"))
$s1.Range("Q1").Value = "total"

# --- Sheet 2 ("建物" / building) ---
# Move the data that lived in row 2 up into row 1 (overwriting the old
# header labels) so it now matches the header layout used on sheet 1,
# and then drop the (now duplicate) row 2.
$s2 = $wb.Worksheets.Item(2)
for ($c = 2; $c -le 8; $c++) {
    $v = $s2.Cells.Item(2, $c).Value()
    $s2.Cells.Item(1, $c).Value = $v
}
$s2.Rows.Item(2).Delete()
